# Update countries & provincias Spain
#
# The source data table ("Pais" sheet) was refreshed:
#  - Some neighboring country rows had their relative order swapped
#    (country name moved to the other row), and
#  - Numbers (Casos totales, Nuevos casos, Casos activos, Recuperados,
#    Casos criticos, Muertes hoy, Muertes) were refreshed for several rows.
#  - The "last updated" timestamp in A1 changed from 08:05 to 09:05.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Mayo de 2020 a las 09:05"

# --- Row 38 / 39: Polonia <-> Ucrania swap + refreshed figures ---
$ws.Range("A38").Value = "Ucrania"
$ws.Range("B38").Value = 23204
$ws.Range("C38").Value = 393
$ws.Range("D38").Value = 9311
$ws.Range("E38").Value = 13197
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 17
$ws.Range("H38").Value = 696

$ws.Range("A39").Value = "Polonia"
$ws.Range("B39").Value = 23155
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 10692
$ws.Range("E39").Value = 11412
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 1051

# --- Row 59 / 60: Bolivia <-> Armenia swap + refreshed figures ---
$ws.Range("A59").Value = "Armenia"
$ws.Range("B59").Value = 8927
$ws.Range("C59").Value = 251
$ws.Range("D59").Value = 3317
$ws.Range("E59").Value = 5483
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 7
$ws.Range("H59").Value = 127

$ws.Range("A60").Value = "Bolivia"
$ws.Range("B60").Value = 8731
$ws.Range("C60").Value = 344
$ws.Range("D60").Value = 749
$ws.Range("E60").Value = 7682
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 7
$ws.Range("H60").Value = 300

# --- Row 75: Hungria refreshed figures ---
$ws.Range("B75").Value = 3867
$ws.Range("C75").Value = 26
$ws.Range("D75").Value = 2142
$ws.Range("E75").Value = 1201
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 7
$ws.Range("H75").Value = 524

# --- Row 114: Letonia refreshed figures ---
$ws.Range("B114").Value = 1065
$ws.Range("C114").Value = 1
$ws.Range("E114").Value = 296

# --- Row 129: Georgia refreshed figures ---
$ws.Range("B129").Value = 757
$ws.Range("C129").Value = 11
$ws.Range("D129").Value = 600
$ws.Range("E129").Value = 145

# --- Row 184: San Martin (Parte Francesa) refreshed figures ---
$ws.Range("B184").Value = 41
$ws.Range("C184").Value = 1
$ws.Range("E184").Value = 5

# --- Row 198 / 199: Curazao <-> Fiyi swap + refreshed figures ---
$ws.Range("A198").Value = "Fiyi"
$ws.Range("D198").Value = 15
$ws.Range("H198").Value = 0

$ws.Range("A199").Value = "Curazao"
$ws.Range("D199").Value = 14
$ws.Range("H199").Value = 1

# --- Row 210 / 211: Montserrat <-> Seychelles swap + refreshed figures ---
$ws.Range("A210").Value = "Seychelles"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# --- Row 213 / 214: Islas Virgenes Britanicas <-> Papua Nueva Guinea swap + refreshed figures ---
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1
